$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the "ODM" (C1) and "Program" (D1) column titles ---
$ws.Cells.Item(1,3).Value = "Program"
$ws.Cells.Item(1,4).Value = "ODM"

# --- New row 25 needs the same style as the other index-column (A) cells ---
$ws.Cells.Item(24,1).Copy()
$ws.Cells.Item(25,1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "SBDPFKBP010T"
$ws.Cells.Item(2,3).Value = "FAIRVIEW"
$ws.Cells.Item(2,4).Value = "PEGATRON"
$ws.Cells.Item(2,5).Value = 5.642222222222221
$ws.Cells.Item(2,6).Value = 202241
$ws.Cells.Item(2,7).Value = 202253
$ws.Cells.Item(2,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(2,9).Value = 5.119069767441861
$ws.Cells.Item(2,10).Value = 195.6622222222222

# Row 3
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "SBDPFKBP020T"
$ws.Cells.Item(3,3).Value = "FAIRVIEW"
$ws.Cells.Item(3,4).Value = "PEGATRON"
$ws.Cells.Item(3,5).Value = 5.653999999999999
$ws.Cells.Item(3,6).Value = 202241
$ws.Cells.Item(3,7).Value = 202253
$ws.Cells.Item(3,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(3,9).Value = 5.738799999999999
$ws.Cells.Item(3,10).Value = 229.552

# Row 4
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "SBDPFKBP512G"
$ws.Cells.Item(4,3).Value = "FAIRVIEW"
$ws.Cells.Item(4,4).Value = "PEGATRON"
$ws.Cells.Item(4,5).Value = 5.642222222222221
$ws.Cells.Item(4,6).Value = 202241
$ws.Cells.Item(4,7).Value = 202253
$ws.Cells.Item(4,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(4,9).Value = 5.270769230769229
$ws.Cells.Item(4,10).Value = 304.5333333333332

# Row 5
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "SBFPF2BU038TES1"
$ws.Cells.Item(5,3).Value = "ADPRRR EE"
$ws.Cells.Item(5,4).Value = "PEGATRON"
$ws.Cells.Item(5,5).Value = ""
$ws.Cells.Item(5,6).Value = 202241
$ws.Cells.Item(5,7).Value = 202253
$ws.Cells.Item(5,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(5,9).Value = ""
$ws.Cells.Item(5,10).Value = ""

# Row 6
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "SBFPF2BU076TES1"
$ws.Cells.Item(6,3).Value = "ADPRRR EE"
$ws.Cells.Item(6,4).Value = "PEGATRON"
$ws.Cells.Item(6,5).Value = ""
$ws.Cells.Item(6,6).Value = 202241
$ws.Cells.Item(6,7).Value = 202253
$ws.Cells.Item(6,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(6,9).Value = ""
$ws.Cells.Item(6,10).Value = ""

# Row 7
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "SBFPF2BU153TES1"
$ws.Cells.Item(7,3).Value = "ADPRRR EE"
$ws.Cells.Item(7,4).Value = "PEGATRON"
$ws.Cells.Item(7,5).Value = ""
$ws.Cells.Item(7,6).Value = 202241
$ws.Cells.Item(7,7).Value = 202253
$ws.Cells.Item(7,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(7,9).Value = ""
$ws.Cells.Item(7,10).Value = ""

# Row 8
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "SBFPF2BV025TES1"
$ws.Cells.Item(8,3).Value = "ADPRRR EE"
$ws.Cells.Item(8,4).Value = "PEGATRON"
$ws.Cells.Item(8,5).Value = 45.39
$ws.Cells.Item(8,6).Value = 202241
$ws.Cells.Item(8,7).Value = 202253
$ws.Cells.Item(8,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(8,9).Value = 45.39
$ws.Cells.Item(8,10).Value = 6354.6

# Row 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "SBFPF2BV025TES1"
$ws.Cells.Item(9,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(9,4).Value = "PEGATRON"
$ws.Cells.Item(9,5).Value = 45.39
$ws.Cells.Item(9,6).Value = 202241
$ws.Cells.Item(9,7).Value = 202253
$ws.Cells.Item(9,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(9,9).Value = 45.39
$ws.Cells.Item(9,10).Value = 1815.6

# Row 10
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "SBFPF2BV076TES1"
$ws.Cells.Item(10,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(10,4).Value = "PEGATRON"
$ws.Cells.Item(10,5).Value = 45.35
$ws.Cells.Item(10,6).Value = 202241
$ws.Cells.Item(10,7).Value = 202253
$ws.Cells.Item(10,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(10,9).Value = 45.35000000000001
$ws.Cells.Item(10,10).Value = 2358.2

# Row 11
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "SBFPF2BV153TES1"
$ws.Cells.Item(11,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(11,4).Value = "PEGATRON"
$ws.Cells.Item(11,5).Value = 45.41
$ws.Cells.Item(11,6).Value = 202241
$ws.Cells.Item(11,7).Value = 202253
$ws.Cells.Item(11,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(11,9).Value = 45.41
$ws.Cells.Item(11,10).Value = 7810.52

# Row 12
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "SBFPF2BV307TES1"
$ws.Cells.Item(12,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(12,4).Value = "PEGATRON"
$ws.Cells.Item(12,5).Value = 120.88
$ws.Cells.Item(12,6).Value = 202241
$ws.Cells.Item(12,7).Value = 202253
$ws.Cells.Item(12,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(12,9).Value = 120.88
$ws.Cells.Item(12,10).Value = 12692.4

# Row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "SBFPFABU038TES1"
$ws.Cells.Item(13,3).Value = "ADPRRR EE"
$ws.Cells.Item(13,4).Value = "PEGATRON"
$ws.Cells.Item(13,5).Value = ""
$ws.Cells.Item(13,6).Value = 202241
$ws.Cells.Item(13,7).Value = 202253
$ws.Cells.Item(13,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(13,9).Value = ""
$ws.Cells.Item(13,10).Value = ""

# Row 14
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "SBFPFABU076TES1"
$ws.Cells.Item(14,3).Value = "ADPRRR EE"
$ws.Cells.Item(14,4).Value = "PEGATRON"
$ws.Cells.Item(14,5).Value = 53.78
$ws.Cells.Item(14,6).Value = 202241
$ws.Cells.Item(14,7).Value = 202253
$ws.Cells.Item(14,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(14,9).Value = 53.78
$ws.Cells.Item(14,10).Value = 1290.72

# Row 15
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "SBFPFABU153TES1"
$ws.Cells.Item(15,3).Value = "ADPRRR EE"
$ws.Cells.Item(15,4).Value = "PEGATRON"
$ws.Cells.Item(15,5).Value = 53.02
$ws.Cells.Item(15,6).Value = 202241
$ws.Cells.Item(15,7).Value = 202253
$ws.Cells.Item(15,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(15,9).Value = 53.02
$ws.Cells.Item(15,10).Value = 1060.4

# Row 16
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "SBFPFABV076TES1"
$ws.Cells.Item(16,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(16,4).Value = "PEGATRON"
$ws.Cells.Item(16,5).Value = 52.96
$ws.Cells.Item(16,6).Value = 202241
$ws.Cells.Item(16,7).Value = 202253
$ws.Cells.Item(16,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(16,9).Value = 52.96
$ws.Cells.Item(16,10).Value = 5084.16

# Row 17
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "SBFPFABV153TES1"
$ws.Cells.Item(17,3).Value = "ADPRRR VE PRQ2"
$ws.Cells.Item(17,4).Value = "PEGATRON"
$ws.Cells.Item(17,5).Value = 53.02
$ws.Cells.Item(17,6).Value = 202241
$ws.Cells.Item(17,7).Value = 202253
$ws.Cells.Item(17,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(17,9).Value = 53.02
$ws.Cells.Item(17,10).Value = 2438.92

# Row 18
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "SBFPFWBV153TES1"
$ws.Cells.Item(18,3).Value = "ADPRRR VE"
$ws.Cells.Item(18,4).Value = "PEGATRON"
$ws.Cells.Item(18,5).Value = 121.16
$ws.Cells.Item(18,6).Value = 202241
$ws.Cells.Item(18,7).Value = 202253
$ws.Cells.Item(18,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(18,9).Value = 121.16
$ws.Cells.Item(18,10).Value = 13085.28

# Row 19
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "SBFPFWBV307TES1"
$ws.Cells.Item(19,3).Value = "ADPRRR VE"
$ws.Cells.Item(19,4).Value = "PEGATRON"
$ws.Cells.Item(19,5).Value = 120.72
$ws.Cells.Item(19,6).Value = 202241
$ws.Cells.Item(19,7).Value = 202253
$ws.Cells.Item(19,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(19,9).Value = 120.72
$ws.Cells.Item(19,10).Value = 12554.88

# Row 20
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "SSDPF2KX012TZES"
$ws.Cells.Item(20,3).Value = "ADPR SE"
$ws.Cells.Item(20,4).Value = "PEGATRON"
$ws.Cells.Item(20,5).Value = 37.7
$ws.Cells.Item(20,6).Value = 202241
$ws.Cells.Item(20,7).Value = 202253
$ws.Cells.Item(20,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(20,9).Value = 37.7
$ws.Cells.Item(20,10).Value = 2865.2

# Row 21
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "SSDPF2KX038TZES"
$ws.Cells.Item(21,3).Value = "ADPR SE"
$ws.Cells.Item(21,4).Value = "PEGATRON"
$ws.Cells.Item(21,5).Value = 37.7
$ws.Cells.Item(21,6).Value = 202241
$ws.Cells.Item(21,7).Value = 202253
$ws.Cells.Item(21,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(21,9).Value = 37.7
$ws.Cells.Item(21,10).Value = 6032

# Row 22
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "SSDPF2KX076TZES"
$ws.Cells.Item(22,3).Value = "ADPR SE"
$ws.Cells.Item(22,4).Value = "PEGATRON"
$ws.Cells.Item(22,5).Value = 37.7
$ws.Cells.Item(22,6).Value = 202241
$ws.Cells.Item(22,7).Value = 202253
$ws.Cells.Item(22,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(22,9).Value = 37.7
$ws.Cells.Item(22,10).Value = 452.4

# Row 23
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "SSDPFINU512GZ1S"
$ws.Cells.Item(23,3).Value = "EH"
$ws.Cells.Item(23,4).Value = "PEGATRON"
$ws.Cells.Item(23,5).Value = 5.79
$ws.Cells.Item(23,6).Value = 202241
$ws.Cells.Item(23,7).Value = 202253
$ws.Cells.Item(23,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(23,9).Value = 5.79
$ws.Cells.Item(23,10).Value = 277.92

# Row 24
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "SSDPFKNU010TZ1S"
$ws.Cells.Item(24,3).Value = "EH"
$ws.Cells.Item(24,4).Value = "PEGATRON"
$ws.Cells.Item(24,5).Value = 5.79
$ws.Cells.Item(24,6).Value = 202241
$ws.Cells.Item(24,7).Value = 202253
$ws.Cells.Item(24,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(24,9).Value = 5.79
$ws.Cells.Item(24,10).Value = 555.84

# Row 25
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "SSDPFKNU512GZ1S"
$ws.Cells.Item(25,3).Value = "EH"
$ws.Cells.Item(25,4).Value = "PEGATRON"
$ws.Cells.Item(25,5).Value = 5.79
$ws.Cells.Item(25,6).Value = 202241
$ws.Cells.Item(25,7).Value = 202253
$ws.Cells.Item(25,8).Value = "ACTIVE, WIP, DONE"
$ws.Cells.Item(25,9).Value = 5.79
$ws.Cells.Item(25,10).Value = 648.48

